$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reposition/resize the workbook window (mirrors the xWindow/windowHeight
#     change in the saved view state). Best-effort: harmless if unsupported. ---
try {
    $win = $excel.ActiveWindow
    $win.Left = 3200
    $win.Height = 16720
} catch {}

# --- Header row 1: re-label column D (old text removed from the shared-string
#     table automatically once nothing else references it); E1/F1 keep their
#     existing text/shared-string slots untouched. ---
$ws.Range("D1").Value = "average incidence rate 2025-50 in >55y at baseline (/100,000/yr)"

# --- Update the incidence-rate figures used by the VE calc (D2:D13); the
#     dependent formulas in column E recalc automatically. ---
$ws.Range("D2:D13").Value = 160.461

# --- Taller header row for the longer wrapped label ---
$ws.Rows(1).RowHeight = 45

# --- Two new explanatory notes below the table ---
$ws.Range("D16").Value = "calc how much due to elderly, times ve times coverage times immunosenesce"
$ws.Range("D17").Value = "risk of react increase- so once wane, reeinter laent so high risk of react"

# --- Move the active selection to reflect where editing left off ---
$ws.Range("D18").Select() | Out-Null
